$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the look of the existing header cells (bold, bordered, centered):
# copy H1's formatting onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# I2:J88 data (87 rows x 2 cols) written as one block via a true 2-D array,
# mirroring how Excel COM automation hands a rectangular Range its .Value.
$data = New-Object 'object[,]' 87,2
$data[0,0] = 7
$data[0,1] = 7
$data[1,0] = 5
$data[1,1] = 6
$data[2,0] = 7
$data[2,1] = 8
$data[3,0] = 6
$data[3,1] = 7
$data[4,0] = 7
$data[4,1] = 7
$data[5,0] = 7
$data[5,1] = 7
$data[6,0] = 10
$data[6,1] = 10
$data[7,0] = 5
$data[7,1] = 5
$data[8,0] = 7
$data[8,1] = 7
$data[9,0] = 8
$data[9,1] = 8
$data[10,0] = 5
$data[10,1] = 6
$data[11,0] = 7
$data[11,1] = 7
$data[12,0] = 6
$data[12,1] = 7
$data[13,0] = 9
$data[13,1] = 9
$data[14,0] = 7
$data[14,1] = 7
$data[15,0] = 8
$data[15,1] = 8
$data[16,0] = 8
$data[16,1] = 8
$data[17,0] = 8
$data[17,1] = 8
$data[18,0] = 7
$data[18,1] = 8
$data[19,0] = 8
$data[19,1] = 8
$data[20,0] = 10
$data[20,1] = 10
$data[21,0] = 7
$data[21,1] = 7
$data[22,0] = 8
$data[22,1] = 8
$data[23,0] = 9
$data[23,1] = 9
$data[24,0] = 9
$data[24,1] = 9
$data[25,0] = 9
$data[25,1] = 9
$data[26,0] = 9
$data[26,1] = 9
$data[27,0] = 9
$data[27,1] = 9
$data[28,0] = 9
$data[28,1] = 9
$data[29,0] = 9
$data[29,1] = 9
$data[30,0] = 9
$data[30,1] = 9
$data[31,0] = 8
$data[31,1] = 9
$data[32,0] = 9
$data[32,1] = 9
$data[33,0] = 8
$data[33,1] = 8
$data[34,0] = 9
$data[34,1] = 9
$data[35,0] = 9
$data[35,1] = 9
$data[36,0] = 9
$data[36,1] = 9
$data[37,0] = 9
$data[37,1] = 10
$data[38,0] = 8
$data[38,1] = 8
$data[39,0] = 9
$data[39,1] = 9
$data[40,0] = 9
$data[40,1] = 9
$data[41,0] = 9
$data[41,1] = 9
$data[42,0] = 8
$data[42,1] = 9
$data[43,0] = 9
$data[43,1] = 9
$data[44,0] = 9
$data[44,1] = 9
$data[45,0] = 8
$data[45,1] = 8
$data[46,0] = 9
$data[46,1] = 9
$data[47,0] = 10
$data[47,1] = 10
$data[48,0] = 9
$data[48,1] = 9
$data[49,0] = 9
$data[49,1] = 9
$data[50,0] = 9
$data[50,1] = 9
$data[51,0] = 9
$data[51,1] = 9
$data[52,0] = 9
$data[52,1] = 9
$data[53,0] = 9
$data[53,1] = 9
$data[54,0] = 9
$data[54,1] = 9
$data[55,0] = 9
$data[55,1] = 9
$data[56,0] = 8
$data[56,1] = 8
$data[57,0] = 9
$data[57,1] = 9
$data[58,0] = 8
$data[58,1] = 9
$data[59,0] = 9
$data[59,1] = 9
$data[60,0] = 9
$data[60,1] = 9
$data[61,0] = 9
$data[61,1] = 9
$data[62,0] = 9
$data[62,1] = 9
$data[63,0] = 9
$data[63,1] = 9
$data[64,0] = 8
$data[64,1] = 9
$data[65,0] = 8
$data[65,1] = 9
$data[66,0] = 8
$data[66,1] = 8
$data[67,0] = 9
$data[67,1] = 9
$data[68,0] = 9
$data[68,1] = 9
$data[69,0] = 9
$data[69,1] = 9
$data[70,0] = 9
$data[70,1] = 9
$data[71,0] = 9
$data[71,1] = 9
$data[72,0] = 9
$data[72,1] = 9
$data[73,0] = 9
$data[73,1] = 9
$data[74,0] = 9
$data[74,1] = 9
$data[75,0] = 9
$data[75,1] = 9
$data[76,0] = 10
$data[76,1] = 11
$data[77,0] = 9
$data[77,1] = 9
$data[78,0] = 9
$data[78,1] = 9
$data[79,0] = 9
$data[79,1] = 9
$data[80,0] = 6
$data[80,1] = 6
$data[81,0] = 5
$data[81,1] = 5
$data[82,0] = 4
$data[82,1] = 4
$data[83,0] = 7
$data[83,1] = 7
$data[84,0] = 4
$data[84,1] = 4
$data[85,0] = 5
$data[85,1] = 5
$data[86,0] = 4
$data[86,1] = 4

$ws.Range("I2:J88").Value = $data

Write-Output "done"
